$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.100.73'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.306.56'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '300.72'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').Value = '97.88'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('E7').Value = '  +3.08%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = '35.84'
$ws.Range('E10').Value = '  -1.38%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = '17.97'
$ws.Range('E13').Value = '  -4.20%  '
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').Value = '2.665.04'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = '2.309.97'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').Value = '43.010.51'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('E19').Value = '  +3.63%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = '6.12'
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('D22').Value = '68.32'
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('D23').Value = '238.17'
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('E24').Value = '  -1.99%  '
$ws.Range('D25').Value = '0.991'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('D28').Value = '25.18'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Value = '166.85'
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('E31').Value = '  -13.52%  '
$ws.Range('E32').Value = '  -5.71%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  +1.45%  '
$ws.Range('D35').Value = '18.28'
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').Value = '0.0690'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('D43').Value = '2.009.15'
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('E45').Value = '  -8.69%  '
$ws.Range('D46').Value = '10.21'
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').Value = '17.43'
$ws.Range('E47').Value = '  -1.52%  '
$ws.Range('E48').Value = '  -2.56%  '
$ws.Range('D49').Value = '54.50'
$ws.Range('E49').Value = '  -2.43%  '
$ws.Range('D50').Value = '2.537.10'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('E51').Value = '  -1.59%  '
